# Updated data types and inconsistent fields with their respective columns:
#  - join_date (E): parse date-like strings into real Excel dates, formatted
#    as "YYYY-MM-DD HH:MM:SS" (also normalizes missing / malformed entries).
#  - performance_score (G): clean up stray "x/1"-style text and fill in blanks
#    with numeric scores.
#  - active (H): convert stray text/blank cells into real booleans.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- join_date (column E) -> real Excel date serials, formatted as date-time ---
$ws.Cells.Item(2, 5).Value = 35886
$ws.Cells.Item(2, 5).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 5).Value = 44575
$ws.Cells.Item(3, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 5).Value = 42660
$ws.Cells.Item(4, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 5).Value = 44792
$ws.Cells.Item(5, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6, 5).Value = 43734
$ws.Cells.Item(6, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7, 5).Value = 44447
$ws.Cells.Item(7, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8, 5).Value = 42265
$ws.Cells.Item(8, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9, 5).Value = 40940
$ws.Cells.Item(9, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 5).Value = 43581
$ws.Cells.Item(10, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 5).Value = 43568
$ws.Cells.Item(11, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12, 5).Value = 44445
$ws.Cells.Item(12, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13, 5).Value = 36495
$ws.Cells.Item(13, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14, 5).Value = 43922
$ws.Cells.Item(14, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value = 44172
$ws.Cells.Item(15, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16, 5).Value = 40940
$ws.Cells.Item(16, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value = 43034
$ws.Cells.Item(17, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 5).Value = 44271
$ws.Cells.Item(18, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(19, 5).Value = 42052
$ws.Cells.Item(19, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(20, 5).Value = 44791
$ws.Cells.Item(20, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21, 5).Value = 44957
$ws.Cells.Item(21, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22, 5).Value = 42474
$ws.Cells.Item(22, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(23, 5).Value = 40940
$ws.Cells.Item(23, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(24, 5).Value = 38867
$ws.Cells.Item(24, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(25, 5).Value = 44797
$ws.Cells.Item(25, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 5).Value = 44204
$ws.Cells.Item(26, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(27, 5).Value = 43097
$ws.Cells.Item(27, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(28, 5).Value = 44630
$ws.Cells.Item(28, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29, 5).Value = 43398
$ws.Cells.Item(29, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(30, 5).Value = 40940
$ws.Cells.Item(30, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(31, 5).Value = 42461
$ws.Cells.Item(31, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(32, 5).Value = 43207
$ws.Cells.Item(32, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(33, 5).Value = 43785
$ws.Cells.Item(33, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 5).Value = 42652
$ws.Cells.Item(34, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 5).Value = 45930
$ws.Cells.Item(35, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 43863
$ws.Cells.Item(36, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(37, 5).Value = 41483
$ws.Cells.Item(37, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(38, 5).Value = 42018
$ws.Cells.Item(38, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(39, 5).Value = 44962
$ws.Cells.Item(39, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 5).Value = 44952
$ws.Cells.Item(40, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(41, 5).Value = 43083
$ws.Cells.Item(41, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(42, 5).Value = 44055
$ws.Cells.Item(42, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(43, 5).Value = 42736
$ws.Cells.Item(43, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(44, 5).Value = 41483
$ws.Cells.Item(44, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(45, 5).Value = 44084
$ws.Cells.Item(45, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 5).Value = 38838
$ws.Cells.Item(46, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(47, 5).Value = 42440
$ws.Cells.Item(47, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 5).Value = 44566
$ws.Cells.Item(48, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(49, 5).Value = 43227
$ws.Cells.Item(49, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(50, 5).Value = 44622
$ws.Cells.Item(50, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(51, 5).Value = 40940
$ws.Cells.Item(51, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(52, 5).Value = 44084
$ws.Cells.Item(52, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(53, 5).Value = 44499
$ws.Cells.Item(53, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 5).Value = 42819
$ws.Cells.Item(54, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(55, 5).Value = 42631
$ws.Cells.Item(55, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 5).Value = 43536
$ws.Cells.Item(56, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(57, 5).Value = 38838
$ws.Cells.Item(57, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(58, 5).Value = 42892
$ws.Cells.Item(58, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(59, 5).Value = 42666
$ws.Cells.Item(59, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(60, 5).Value = 44214
$ws.Cells.Item(60, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(61, 5).Value = 44177
$ws.Cells.Item(61, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- performance_score (column G) -> numeric values ---
$ws.Cells.Item(5, 7).Value = 0.33
$ws.Cells.Item(8, 7).Value = 0.68
$ws.Cells.Item(36, 7).Value = 0.8100000000000001
$ws.Cells.Item(43, 7).Value = 0.68
$ws.Cells.Item(45, 7).Value = 0.43
$ws.Cells.Item(48, 7).Value = 0.01
$ws.Cells.Item(57, 7).Value = 0.07000000000000001

# --- active (column H) -> boolean values ---
$ws.Cells.Item(9, 8).Value = $true
$ws.Cells.Item(19, 8).Value = $true
$ws.Cells.Item(22, 8).Value = $true
$ws.Cells.Item(41, 8).Value = $true
$ws.Cells.Item(48, 8).Value = $true
$ws.Cells.Item(49, 8).Value = $true
$ws.Cells.Item(54, 8).Value = $true
$ws.Cells.Item(56, 8).Value = $false
$ws.Cells.Item(57, 8).Value = $false
$ws.Cells.Item(60, 8).Value = $false
